# Auto-generated edit script: refreshes cached market-board price/profit
# columns (H-N) across all eight crafting-leve sheets with updated values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 249.77777
$ws.Range("I2").Value = 284.83334
$ws.Range("J2").Value = 179.66667
$ws.Range("K2").Value = 284.83334
$ws.Range("L2").Value = 179.66667
$ws.Range("M2").Value = -171.83334
$ws.Range("N2").Value = -405.66667
$ws.Range("H9").Value = 5699.222
$ws.Range("I9").Value = 7855.385
$ws.Range("J9").Value = 93.2
$ws.Range("K9").Value = 7855.385
$ws.Range("L9").Value = 93.2
$ws.Range("M9").Value = -7686.385
$ws.Range("N9").Value = -431.2
$ws.Range("H33").Value = 301.14285
$ws.Range("I33").Value = 284.66666
$ws.Range("K33").Value = 284.66666
$ws.Range("M33").Value = -55.66665999999998
$ws.Range("H40").Value = 4918.129
$ws.Range("I40").Value = 2995
$ws.Range("J40").Value = 6501.8823
$ws.Range("K40").Value = 2995
$ws.Range("L40").Value = 6501.8823
$ws.Range("M40").Value = -2820
$ws.Range("N40").Value = -6851.8823
$ws.Range("H41").Value = 1053.0834
$ws.Range("I41").Value = 1366.8889
$ws.Range("J41").Value = 111.666664
$ws.Range("K41").Value = 1366.8889
$ws.Range("L41").Value = 111.666664
$ws.Range("M41").Value = -926.8888999999999
$ws.Range("N41").Value = -991.666664
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H62").Value = 1248.6666
$ws.Range("I62").Value = 1219.8572
$ws.Range("K62").Value = 1219.8572
$ws.Range("M62").Value = -595.8571999999999
$ws.Range("H65").Value = 1248.6666
$ws.Range("I65").Value = 1219.8572
$ws.Range("K65").Value = 6099.286
$ws.Range("M65").Value = -2979.286
$ws.Range("H70").Value = 3442.2646
$ws.Range("I70").Value = 1110.8948
$ws.Range("J70").Value = 6395.3335
$ws.Range("K70").Value = 3332.6844
$ws.Range("L70").Value = 19186.0005
$ws.Range("M70").Value = -3062.6844
$ws.Range("N70").Value = -19726.0005
$ws.Range("H73").Value = 3442.2646
$ws.Range("I73").Value = 1110.8948
$ws.Range("J73").Value = 6395.3335
$ws.Range("K73").Value = 3332.6844
$ws.Range("L73").Value = 19186.0005
$ws.Range("M73").Value = -2396.6844
$ws.Range("N73").Value = -21058.0005
$ws.Range("H82").Value = 3580.5789
$ws.Range("I82").Value = 1884.2354
$ws.Range("J82").Value = 17999.5
$ws.Range("K82").Value = 5652.706200000001
$ws.Range("L82").Value = 53998.5
$ws.Range("M82").Value = -5246.706200000001
$ws.Range("N82").Value = -54810.5
$ws.Range("H85").Value = 3580.5789
$ws.Range("I85").Value = 1884.2354
$ws.Range("J85").Value = 17999.5
$ws.Range("K85").Value = 5652.706200000001
$ws.Range("L85").Value = 53998.5
$ws.Range("M85").Value = -4248.706200000001
$ws.Range("N85").Value = -56806.5
$ws.Range("H86").Value = 5830.1113
$ws.Range("J86").Value = 2133.6667
$ws.Range("L86").Value = 2133.6667
$ws.Range("N86").Value = -4379.6667
$ws.Range("I88").Value = 333333340
$ws.Range("J88").Value = 7158456
$ws.Range("K88").Value = 333333340
$ws.Range("L88").Value = 7158456
$ws.Range("M88").Value = -333332934
$ws.Range("N88").Value = -7159268
$ws.Range("H89").Value = 5830.1113
$ws.Range("J89").Value = 2133.6667
$ws.Range("L89").Value = 10668.3335
$ws.Range("N89").Value = -21900.3335
$ws.Range("I91").Value = 333333340
$ws.Range("J91").Value = 7158456
$ws.Range("K91").Value = 333333340
$ws.Range("L91").Value = 7158456
$ws.Range("M91").Value = -333331936
$ws.Range("N91").Value = -7161264
$ws.Range("H135").Value = 31250348
$ws.Range("I135").Value = 31250348
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 281253132
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -281250597
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2344.3
$ws.Range("I137").Value = 1551.75
$ws.Range("J137").Value = 3533.125
$ws.Range("K137").Value = 4655.25
$ws.Range("L137").Value = 10599.375
$ws.Range("M137").Value = -2105.25
$ws.Range("N137").Value = -15699.375
$ws.Range("H138").Value = 3466.65
$ws.Range("I138").Value = 2520.4
$ws.Range("K138").Value = 7561.200000000001
$ws.Range("M138").Value = -2421.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1633.4
$ws.Range("J5").Value = 3540
$ws.Range("L5").Value = 3540
$ws.Range("N5").Value = -3764
$ws.Range("H27").Value = 5000
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5368
$ws.Range("H31").Value = 15500
$ws.Range("I31").Value = 15500
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 15500
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -15206
$ws.Range("N31").ClearContents()
$ws.Range("H61").Value = 62502564
$ws.Range("I61").Value = 111112450
$ws.Range("J61").Value = 4142.143
$ws.Range("K61").Value = 111112450
$ws.Range("L61").Value = 4142.143
$ws.Range("M61").Value = -111112238
$ws.Range("N61").Value = -4566.143
$ws.Range("H74").Value = 43480770
$ws.Range("I74").Value = 47620844
$ws.Range("K74").Value = 47620844
$ws.Range("M74").Value = -47619970
$ws.Range("H77").Value = 43480770
$ws.Range("I77").Value = 47620844
$ws.Range("K77").Value = 238104220
$ws.Range("M77").Value = -238099852
$ws.Range("H97").Value = 622.6667
$ws.Range("I97").Value = 793.63635
$ws.Range("J97").Value = 280.72726
$ws.Range("K97").Value = 793.63635
$ws.Range("L97").Value = 280.72726
$ws.Range("M97").Value = -297.63635
$ws.Range("N97").Value = -1272.72726
$ws.Range("H122").Value = 4957.3184
$ws.Range("I122").Value = 5362.273
$ws.Range("J122").Value = 4552.364
$ws.Range("K122").Value = 16086.819
$ws.Range("L122").Value = 13657.092
$ws.Range("M122").Value = -13636.819
$ws.Range("N122").Value = -18557.092
$ws.Range("H132").Value = 9094534
$ws.Range("I132").Value = 12502172
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 37506516
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -37503986
$ws.Range("N132").Value = -27560
$ws.Range("H136").Value = 62502564
$ws.Range("I136").Value = 111112450
$ws.Range("J136").Value = 4142.143
$ws.Range("K136").Value = 333337350
$ws.Range("L136").Value = 12426.429
$ws.Range("M136").Value = -333334800
$ws.Range("N136").Value = -17526.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1633.4
$ws.Range("J4").Value = 3540
$ws.Range("L4").Value = 3540
$ws.Range("N4").Value = -3770
$ws.Range("H21").Value = 25500
$ws.Range("J21").Value = 25500
$ws.Range("L21").Value = 25500
$ws.Range("N21").Value = -25972
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H54").Value = 40000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H75").Value = 12237.75
$ws.Range("I75").Value = 9039
$ws.Range("J75").Value = 17569
$ws.Range("K75").Value = 9039
$ws.Range("L75").Value = 17569
$ws.Range("M75").Value = -8103
$ws.Range("N75").Value = -19441
$ws.Range("H78").Value = 12237.75
$ws.Range("I78").Value = 9039
$ws.Range("J78").Value = 17569
$ws.Range("K78").Value = 27117
$ws.Range("L78").Value = 52707
$ws.Range("M78").Value = -22437
$ws.Range("N78").Value = -62067
$ws.Range("H82").Value = 9668
$ws.Range("I82").Value = 9668
$ws.Range("K82").Value = 9668
$ws.Range("M82").Value = -9285
$ws.Range("H85").Value = 9668
$ws.Range("I85").Value = 9668
$ws.Range("K85").Value = 9668
$ws.Range("M85").Value = -8342
$ws.Range("H86").Value = 2936.8
$ws.Range("I86").Value = 3152
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 3152
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -2029
$ws.Range("N86").Value = -3246
$ws.Range("H89").Value = 2936.8
$ws.Range("I89").Value = 3152
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 15760
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -10144
$ws.Range("N89").Value = -16232
$ws.Range("H94").Value = 3006.5334
$ws.Range("I94").Value = 3006.5334
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3006.5334
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2555.5334
$ws.Range("N94").ClearContents()
$ws.Range("H99").Value = 2199.8462
$ws.Range("I99").Value = 2063
$ws.Range("J99").Value = 2240.9
$ws.Range("K99").Value = 2063
$ws.Range("L99").Value = 2240.9
$ws.Range("M99").Value = -565
$ws.Range("N99").Value = -5236.9
$ws.Range("H103").Value = 24330
$ws.Range("J103").Value = 24330
$ws.Range("L103").Value = 24330
$ws.Range("N103").Value = -26674
$ws.Range("H117").Value = 32993
$ws.Range("J117").Value = 32993
$ws.Range("L117").Value = 32993
$ws.Range("N117").Value = -42171
$ws.Range("H134").Value = 19232528
$ws.Range("I134").Value = 20834968
$ws.Range("K134").Value = 62504904
$ws.Range("M134").Value = -62502369

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6667336.5
$ws.Range("I4").Value = 10
$ws.Range("K4").Value = 10
$ws.Range("M4").Value = 102
$ws.Range("H16").Value = 1372618.4
$ws.Range("I16").Value = 1568135.2
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 1568135.2
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -1567848.2
$ws.Range("N16").Value = -4574
$ws.Range("H43").Value = 12999.2
$ws.Range("J43").Value = 12999.2
$ws.Range("L43").Value = 12999.2
$ws.Range("N43").Value = -13367.2
$ws.Range("H51").Value = 44996
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 44996
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 44996
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -46468
$ws.Range("H58").Value = 13518088
$ws.Range("I58").Value = 23815456
$ws.Range("J58").Value = 2792.3125
$ws.Range("K58").Value = 23815456
$ws.Range("L58").Value = 2792.3125
$ws.Range("M58").Value = -23815253
$ws.Range("N58").Value = -3198.3125
$ws.Range("H60").Value = 18466.334
$ws.Range("I60").Value = 7700
$ws.Range("J60").Value = 39999
$ws.Range("K60").Value = 7700
$ws.Range("L60").Value = 39999
$ws.Range("M60").Value = -7189
$ws.Range("N60").Value = -41021
$ws.Range("H61").Value = 44996
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 44996
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 44996
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -45692
$ws.Range("H62").Value = 4968.2
$ws.Range("I62").Value = 4968.2
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4968.2
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4344.2
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4968.2
$ws.Range("I65").Value = 4968.2
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 24841
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -21721
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 47483.75
$ws.Range("I74").Value = 40000
$ws.Range("J74").Value = 49978.332
$ws.Range("K74").Value = 40000
$ws.Range("L74").Value = 49978.332
$ws.Range("M74").Value = -39126
$ws.Range("N74").Value = -51726.332
$ws.Range("H77").Value = 47483.75
$ws.Range("I77").Value = 40000
$ws.Range("J77").Value = 49978.332
$ws.Range("K77").Value = 120000
$ws.Range("L77").Value = 149934.996
$ws.Range("M77").Value = -115632
$ws.Range("N77").Value = -158670.996
$ws.Range("H86").Value = 12889.571
$ws.Range("J86").Value = 15549.875
$ws.Range("L86").Value = 15549.875
$ws.Range("N86").Value = -17795.875
$ws.Range("H89").Value = 12889.571
$ws.Range("J89").Value = 15549.875
$ws.Range("L89").Value = 77749.375
$ws.Range("N89").Value = -88981.375
$ws.Range("H101").Value = 12999.2
$ws.Range("J101").Value = 12999.2
$ws.Range("L101").Value = 12999.2
$ws.Range("N101").Value = -19489.2
$ws.Range("H107").Value = 367763.25
$ws.Range("I107").Value = 505628.1
$ws.Range("K107").Value = 505628.1
$ws.Range("M107").Value = -503708.1
$ws.Range("H109").Value = 30997.5
$ws.Range("J109").Value = 30997.5
$ws.Range("L109").Value = 30997.5
$ws.Range("N109").Value = -33077.5
$ws.Range("H113").Value = 1372618.4
$ws.Range("I113").Value = 1568135.2
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 1568135.2
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -1565965.2
$ws.Range("N113").Value = -8340
$ws.Range("H122").Value = 2106
$ws.Range("I122").Value = 2156.8572
$ws.Range("J122").Value = 1394
$ws.Range("K122").Value = 6470.571599999999
$ws.Range("L122").Value = 4182
$ws.Range("M122").Value = -4020.571599999999
$ws.Range("N122").Value = -9082
$ws.Range("H132").Value = 166669330
$ws.Range("I132").Value = 500001000
$ws.Range("J132").Value = 3497.25
$ws.Range("K132").Value = 1500003000
$ws.Range("L132").Value = 10491.75
$ws.Range("M132").Value = -1500000470
$ws.Range("N132").Value = -15551.75
$ws.Range("H134").Value = 17931516
$ws.Range("I134").Value = 19310634
$ws.Range("K134").Value = 57931902
$ws.Range("M134").Value = -57929367
$ws.Range("H136").Value = 13518088
$ws.Range("I136").Value = 23815456
$ws.Range("J136").Value = 2792.3125
$ws.Range("K136").Value = 71446368
$ws.Range("L136").Value = 8376.9375
$ws.Range("M136").Value = -71443818
$ws.Range("N136").Value = -13476.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 40
$ws.Range("J10").Value = 40
$ws.Range("L10").Value = 120
$ws.Range("N10").Value = -398
$ws.Range("H12").Value = 313.4737
$ws.Range("I12").Value = 45.333332
$ws.Range("K12").Value = 135.999996
$ws.Range("M12").Value = 37.00000399999999
$ws.Range("H33").Value = 960.7143
$ws.Range("I33").Value = 559.6
$ws.Range("J33").Value = 1183.5555
$ws.Range("K33").Value = 3357.6
$ws.Range("L33").Value = 7101.333
$ws.Range("M33").Value = -3074.6
$ws.Range("N33").Value = -7667.333
$ws.Range("H92").Value = 50
$ws.Range("J92").Value = 50
$ws.Range("L92").Value = 150
$ws.Range("N92").Value = -2646
$ws.Range("H122").Value = 1366.6666
$ws.Range("J122").Value = 1425
$ws.Range("L122").Value = 12825
$ws.Range("N122").Value = -17725
$ws.Range("H131").Value = 1765.4138
$ws.Range("I131").Value = 1364.3636
$ws.Range("J131").Value = 2010.5
$ws.Range("K131").Value = 4093.0908
$ws.Range("L131").Value = 6031.5
$ws.Range("M131").Value = 946.9092000000001
$ws.Range("N131").Value = -16111.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 193.2
$ws.Range("I2").Value = 170.8
$ws.Range("J2").Value = 215.6
$ws.Range("K2").Value = 170.8
$ws.Range("L2").Value = 215.6
$ws.Range("M2").Value = -57.80000000000001
$ws.Range("N2").Value = -441.6
$ws.Range("H97").Value = 3429.4
$ws.Range("I97").Value = 3512.375
$ws.Range("K97").Value = 3512.375
$ws.Range("M97").Value = -3016.375
$ws.Range("H98").Value = 114155.4
$ws.Range("J98").Value = 114155.4
$ws.Range("L98").Value = 114155.4
$ws.Range("N98").Value = -120145.4
$ws.Range("H102").Value = 6122.7
$ws.Range("I102").Value = 5174.875
$ws.Range("K102").Value = 5174.875
$ws.Range("M102").Value = -3552.875
$ws.Range("H122").Value = 93966.92999999999
$ws.Range("I122").Value = 114412.45
$ws.Range("J122").Value = 19000
$ws.Range("K122").Value = 343237.35
$ws.Range("L122").Value = 57000
$ws.Range("M122").Value = -340787.35
$ws.Range("N122").Value = -61900
$ws.Range("H132").Value = 8931104
$ws.Range("I132").Value = 10419138
$ws.Range("J132").Value = 2896.5
$ws.Range("K132").Value = 31257414
$ws.Range("L132").Value = 8689.5
$ws.Range("M132").Value = -31254884
$ws.Range("N132").Value = -13749.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1936.2222
$ws.Range("I16").Value = 2144.8572
$ws.Range("J16").Value = 1206
$ws.Range("K16").Value = 2144.8572
$ws.Range("L16").Value = 1206
$ws.Range("M16").Value = -1974.8572
$ws.Range("N16").Value = -1546
$ws.Range("H32").Value = 6000
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H40").Value = 4182.6665
$ws.Range("I40").Value = 4019.4
$ws.Range("K40").Value = 4019.4
$ws.Range("M40").Value = -3883.4
$ws.Range("H55").Value = 498.80646
$ws.Range("I55").Value = 313.11765
$ws.Range("J55").Value = 724.2857
$ws.Range("K55").Value = 313.11765
$ws.Range("L55").Value = 724.2857
$ws.Range("M55").Value = -140.11765
$ws.Range("N55").Value = -1070.2857
$ws.Range("H82").Value = 1546.3871
$ws.Range("I82").Value = 1551.5217
$ws.Range("J82").Value = 1531.625
$ws.Range("K82").Value = 1551.5217
$ws.Range("L82").Value = 1531.625
$ws.Range("M82").Value = -1190.5217
$ws.Range("N82").Value = -2253.625
$ws.Range("H85").Value = 1546.3871
$ws.Range("I85").Value = 1551.5217
$ws.Range("J85").Value = 1531.625
$ws.Range("K85").Value = 1551.5217
$ws.Range("L85").Value = 1531.625
$ws.Range("M85").Value = -303.5217
$ws.Range("N85").Value = -4027.625
$ws.Range("H87").Value = 10526
$ws.Range("J87").Value = 10526
$ws.Range("L87").Value = 10526
$ws.Range("N87").Value = -12772
$ws.Range("H90").Value = 10526
$ws.Range("J90").Value = 10526
$ws.Range("L90").Value = 31578
$ws.Range("N90").Value = -42810
$ws.Range("H93").Value = 2014.6471
$ws.Range("I93").Value = 1054.5385
$ws.Range("J93").Value = 5135
$ws.Range("K93").Value = 1054.5385
$ws.Range("L93").Value = 5135
$ws.Range("M93").Value = 193.4614999999999
$ws.Range("N93").Value = -7631
$ws.Range("H96").Value = 42000
$ws.Range("J96").Value = 42000
$ws.Range("L96").Value = 42000
$ws.Range("N96").Value = -47492
$ws.Range("H122").Value = 12226.5
$ws.Range("I122").Value = 11115.1
$ws.Range("K122").Value = 33345.3
$ws.Range("M122").Value = -30895.3
$ws.Range("H132").Value = 18464836
$ws.Range("I132").Value = 19203148
$ws.Range("K132").Value = 57609444
$ws.Range("M132").Value = -57606914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28108
$ws.Range("I41").Value = 5770.5
$ws.Range("J41").Value = 42999.668
$ws.Range("K41").Value = 5770.5
$ws.Range("L41").Value = 42999.668
$ws.Range("M41").Value = -5380.5
$ws.Range("N41").Value = -43779.668
$ws.Range("H96").Value = 3093
$ws.Range("I96").Value = 1263.3636
$ws.Range("K96").Value = 1263.3636
$ws.Range("M96").Value = 109.6364000000001
$ws.Range("H101").Value = 189995
$ws.Range("J101").Value = 189995
$ws.Range("L101").Value = 189995
$ws.Range("N101").Value = -196485
$ws.Range("H122").Value = 1683.1666
$ws.Range("I122").Value = 1683.1666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5049.4998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2599.4998
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 4583
$ws.Range("I126").Value = 4583
$ws.Range("K126").Value = 13749
$ws.Range("M126").Value = -11279
$ws.Range("H132").Value = 8337984
$ws.Range("I132").Value = 10872809
$ws.Range("J132").Value = 9273
$ws.Range("K132").Value = 32618427
$ws.Range("L132").Value = 27819
$ws.Range("M132").Value = -32615897
$ws.Range("N132").Value = -32879
